# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns for both the zh-cn and de-de handback rows, links the
# newly-populated target-file cells back to the source doc, and updates the
# Overview sheet's per-language status text now that the round trip is done.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec7b46633fa37fd9f13edc2785345ecfe1537d70/e2e/"

# ---------------------------------------------------------------------------
# Overview sheet: both language columns move from "Ready for handoff" to
# "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet: target file + handback file + handback datetime for both rows.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2:C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("I2").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($repoBase + "a.md"), "", "", "a.md")
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 20:46:40"

$zhcn.Range("I3").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($repoBase + "a.md"), "", "", "a.md")
$zhcn.Range("I3").Font.Underline = $true
$zhcn.Range("I3").Font.Color = 15570276
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-31 20:46:40"

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet: target file + handback file + handback datetime for both rows.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2:C3").Value = "Handed back: in sync with en-US"

$dede.Range("I2").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("I2"), ($repoBase + "a.md"), "", "", "a.md")
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = 15570276
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 20:46:47"

$dede.Range("I3").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("I3"), ($repoBase + "a.md"), "", "", "a.md")
$dede.Range("I3").Font.Underline = $true
$dede.Range("I3").Font.Color = 15570276
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-31 20:46:47"

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(10).ColumnWidth = 40
